$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (player name, position, team) for rows 2-18
$data = @(
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls"),
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic")
)

# Clear out the previous data range (rows 2 through 19) before rewriting it
$ws.Range("A2:C19").ClearContents()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
